$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Block 1 (rows 4-8): fill in columns D and H ---
$ws.Range("D4").Value = 31
$ws.Range("H4").Value = 83
$ws.Range("D5").Value = 66
$ws.Range("H5").Value = 84
$ws.Range("D6").Value = 19
$ws.Range("H6").Value = 22
$ws.Range("D7").Value = 17
$ws.Range("H7").Value = 22
$ws.Range("D8").Value = 18
$ws.Range("H8").Value = 21

# --- Block 2 (rows 12-16): fill in columns D and H ---
$ws.Range("D12").Value = 17
$ws.Range("H12").Value = 79
$ws.Range("D13").Value = 19
$ws.Range("H13").Value = 18
$ws.Range("D14").Value = 44
$ws.Range("H14").Value = 16
$ws.Range("D15").Value = 44
$ws.Range("H15").Value = 18
$ws.Range("D16").Value = 51
$ws.Range("H16").Value = 16

# --- Block 3 (rows 20-24): fill in columns D and H ---
$ws.Range("D20").Value = 17
$ws.Range("H20").Value = 19
$ws.Range("D21").Value = 76
$ws.Range("H21").Value = 16
$ws.Range("D22").Value = 78
$ws.Range("H22").Value = 47
$ws.Range("D23").Value = 22
$ws.Range("H23").Value = 418
$ws.Range("D24").Value = 16
$ws.Range("H24").Value = 67

# --- New trailing row 25: blank placeholder cell under the D column,
# mirroring the spacer cell pattern used elsewhere in the sheet (e.g. D9) ---
$ws.Range("D25").WrapText = $False

# --- Selection moves to the last-edited cell ---
$ws.Range("H24").Select()
